$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add hyperlink + text to F3 (email address), keeping the existing "Hyperlink" style (s=1)
$ws.Range("F3").Value = "rrhh@galqui.com"
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:rrhh@galqui.com") | Out-Null
$ws.Range("F3").Style = "Hipervínculo"

# Column J: add a value of 3 for every data row (3-13)
for ($r = 3; $r -le 13; $r++) {
    $ws.Cells.Item($r, 10).Value = 3
}

# Update column I values for rows 7, 9, 11 from 1 to 2
$ws.Range("I7").Value = 2
$ws.Range("I9").Value = 2
$ws.Range("I11").Value = 2

# Update the view: scroll so column B is the top-left visible column, and select J3:J13 with active cell J3
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("B1").Column
$ws.Range("J3:J13").Select()
